$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A..D) for rows 2..10
$varLabel = "Diferença 2025/04 - 2024/04"

$data = @(
    @("Amapá", $varLabel, 2.08, "1º"),
    @("Bahia", $varLabel, 1.88, "2º"),
    @("Rio Grande do Norte", $varLabel, 1.74, "3º"),
    @("Paraíba", $varLabel, 1.71, "4º"),
    @("Rio de Janeiro", $varLabel, 1.55, "5º"),
    @("Rio Grande do Sul", $varLabel, 1.52, "6º"),
    @("Sergipe", $varLabel, 1.18, "11º"),
    @("Brasil", $varLabel, 1.13, ""),
    @("Nordeste", $varLabel, 1.2, "")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
